$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (114 and 115) after the existing data (last row 113)
$ws.Range("A114").Value = 113
$ws.Range("B114").Value = 1
$ws.Range("C114").Value = "2024-06-17 08:16:48"
$ws.Range("D114").Value = 200
$ws.Range("E114").Value = 8

$ws.Range("A115").Value = 114
$ws.Range("B115").Value = 2
$ws.Range("C115").Value = "2024-06-17 08:16:48"
$ws.Range("D115").Value = 200
$ws.Range("E115").Value = 0
